# Preserve whitespace in cells.
# Adds a new worksheet "Test Whitespace" (after the existing "Cell Values"
# sheet) containing a single cell A1 whose value is a whitespace-only
# string ("    ", four spaces) - this exercises / demonstrates that
# whitespace-only string values are preserved (round-tripped with
# xml:space="preserve") rather than being trimmed away.

$wb = $excel.ActiveWorkbook

# Remember the sheet that is active before we start, so we can restore the
# original selection/active-tab once the new sheet has been added.
$originalActiveSheetName = $wb.ActiveSheet.Name

# Add the new worksheet right after the last existing sheet so it lands at
# the end of the tab strip (matching the workbook.xml <x:sheets> order in
# the diff: "Cell Values" stays first, "Test Whitespace" is appended).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Test Whitespace"

# A1 holds a whitespace-only string value - must round-trip untouched.
$newSheet.Range("A1").Value = "    "

# Restore the original active sheet/tab selection.
$wb.Worksheets.Item($originalActiveSheetName).Activate()
